# Update the cryptocurrency price/volume figures per the Nov 2 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$cellRef, [string]$newValue) {
    $rng = $ws.Range($cellRef)
    # Force Text format first so numeric-looking strings (e.g. "11.33",
    # "35.534.74") are stored as literal text, matching the source data,
    # then restore the default "Normal" style so no stray number format
    # is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextCell "D2" "35.534.74"
Set-TextCell "E2" "  +3.14%  "
Set-TextCell "D3" "1.844.70"
Set-TextCell "E3" "  +2.51%  "
Set-TextCell "E4" "  +0.26%  "
Set-TextCell "D5" "232.50"
Set-TextCell "E5" "  +3.63%  "
Set-TextCell "D6" "0.620"
Set-TextCell "E6" "  +3.20%  "
Set-TextCell "E7" "  +0.21%  "
Set-TextCell "D8" "44.15"
Set-TextCell "E8" "  +13.11%  "
Set-TextCell "E9" "  +8.35%  "
Set-TextCell "D10" "0.0699"
Set-TextCell "E10" "  +4.76%  "
Set-TextCell "E11" "  +2.87%  "
Set-TextCell "E12" "  +2.46%  "
Set-TextCell "B13" "Chainlink"
Set-TextCell "C13" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D13" "11.33"
Set-TextCell "E13" "  +4.09%  "
Set-TextCell "B14" "WrappedEther"
Set-TextCell "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D14" "1.835.01"
Set-TextCell "E14" "  +1.89%  "
Set-TextCell "D15" "0.675"
Set-TextCell "E15" "  +7.72%  "
Set-TextCell "D16" "4.73"
Set-TextCell "E16" "  +8.65%  "
Set-TextCell "D17" "35.516.67"
Set-TextCell "E17" "  +3.22%  "
Set-TextCell "D18" "70.55"
Set-TextCell "E18" "  +3.85%  "
Set-TextCell "D19" "0.0₃0803"
Set-TextCell "E19" "  +5.20%  "
Set-TextCell "D20" "244.39"
Set-TextCell "E20" "  +2.34%  "
Set-TextCell "D21" "12.08"
Set-TextCell "E21" "  +9.14%  "
Set-TextCell "E22" "  +13.85%  "
Set-TextCell "E23" "  +0.23%  "
Set-TextCell "D24" "2.26"
Set-TextCell "E24" "  +5.08%  "
Set-TextCell "D25" "171.06"
Set-TextCell "E25" "  +0.42%  "
Set-TextCell "D26" "8.02"
Set-TextCell "E26" "  +4.86%  "
Set-TextCell "D27" "17.84"
Set-TextCell "E27" "  +1.60%  "
Set-TextCell "E28" "  +1.20%  "
Set-TextCell "D29" "1.58"
Set-TextCell "E29" "  +29.40%  "
Set-TextCell "E30" "  +0.23%  "
Set-TextCell "D31" "3.331.01"
Set-TextCell "E31" "  +37.10%  "
Set-TextCell "D32" "0.0554"
Set-TextCell "E32" "  +8.22%  "
Set-TextCell "D33" "4.12"
Set-TextCell "E33" "  +7.45%  "
Set-TextCell "D34" "3.94"
Set-TextCell "E34" "  +5.55%  "
Set-TextCell "D35" "1.84"
Set-TextCell "E35" "  +1.87%  "
Set-TextCell "D36" "95.29"
Set-TextCell "E36" "  +16.88%  "
Set-TextCell "D37" "0.693"
Set-TextCell "E37" "  +8.51%  "
Set-TextCell "E38" "  +8.01%  "
Set-TextCell "D39" "1.349.92"
Set-TextCell "E39" "  +3.65%  "
Set-TextCell "E40" "  +5.96%  "
Set-TextCell "D41" "2.44"
Set-TextCell "E41" "  +6.19%  "
Set-TextCell "D42" "15.34"
Set-TextCell "E42" "  +9.68%  "
Set-TextCell "D43" "1.02"
Set-TextCell "E43" "  +7.99%  "
Set-TextCell "E44" "  +3.50%  "
Set-TextCell "E45" "  +0.90%  "
Set-TextCell "E46" "  +0.28%  "
Set-TextCell "E47" "  +9.74%  "
Set-TextCell "D48" "0.0520"
Set-TextCell "E48" "  +0.39%  "
Set-TextCell "E49" "  +2.89%  "
Set-TextCell "E50" "  +0.28%  "
Set-TextCell "D51" "102.58"
Set-TextCell "E51" "  +0.87%  "
